# Update attendee counts ("想去人数") and a couple of min ticket prices
# ("最低票价") across the four sheets of the workbook, reflecting the
# regenerated site data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 2793
$ws.Range("F10").Value = 366
$ws.Range("F11").Value = 18
$ws.Range("F12").Value = 302
$ws.Range("F14").Value = 5852
$ws.Range("F16").Value = 1033
$ws.Range("F17").Value = 86
$ws.Range("G17").Value = 29.9
$ws.Range("F19").Value = 86
$ws.Range("F21").Value = 11
$ws.Range("F22").Value = 1277
$ws.Range("F25").Value = 2037
$ws.Range("F26").Value = 155
$ws.Range("F29").Value = 3229

# --- Sheet "演出" (performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 25
$ws.Range("F17").Value = 57
$ws.Range("F24").Value = 4042
$ws.Range("F26").Value = 10
$ws.Range("F32").Value = 9

# --- Sheet "本地生活" (local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 1463
$ws.Range("F9").Value = 410
$ws.Range("F12").Value = 604

# --- Sheet "全部类型" (all types / combined) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 1463
$ws.Range("F8").Value = 410
$ws.Range("F13").Value = 25
$ws.Range("F14").Value = 2793
$ws.Range("F17").Value = 604
$ws.Range("F18").Value = 366
$ws.Range("F21").Value = 302
$ws.Range("F23").Value = 5852
$ws.Range("F25").Value = 1033
$ws.Range("F27").Value = 86
$ws.Range("G27").Value = 29.9
$ws.Range("F29").Value = 86
$ws.Range("F33").Value = 57
$ws.Range("F38").Value = 1277
$ws.Range("F43").Value = 2037
$ws.Range("F46").Value = 155
$ws.Range("F49").Value = 3229
